$d = $word.ActiveDocument

# Locate the end of the final content paragraph ("...istent with our
# training data.") so the new paragraphs land right after it, before the
# trailing <w:bookmarkStart/bookmarkEnd name="_GoBack"/> that currently
# sits at the very end of that paragraph.
$anchor = $d.Content
$found = $anchor.Find.Execute("istent with our training data.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find anchor text 'istent with our training data.'"
}
$insertPoint = $d.Range($anchor.End, $anchor.End)

# Build the nine new paragraphs as raw WordprocessingML so the run
# boundaries come through untouched (no auto-coalescing of adjacent
# same-formatted runs) and the "empty" paragraphs stay truly empty
# (<w:p/>, no stray run). The relocated bookmark is (re)created here with
# a scratch id of 1; the stale copy still wrapping the original last
# paragraph is removed afterwards and Word renumbers the survivor back to
# id 0 on save.
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$newParagraphsXml = @"
<w:p $wNs/>
<w:p $wNs/>
<w:p $wNs><w:r><w:t xml:space="preserve">Cross v with </w:t></w:r><w:r><w:t xml:space="preserve">new </w:t></w:r><w:r><w:t>features</w:t></w:r></w:p>
<w:p $wNs><w:r><w:t>1.054171752575413</w:t></w:r></w:p>
<w:p $wNs/>
<w:p $wNs><w:r><w:t>Cross v without</w:t></w:r><w:r><w:t xml:space="preserve"> new features</w:t></w:r><w:bookmarkStart w:id="1" w:name="_GoBack"/><w:bookmarkEnd w:id="1"/></w:p>
<w:p $wNs><w:r><w:t>1.0612860270645978</w:t></w:r></w:p>
<w:p $wNs/>
<w:p $wNs/>
"@
[void]$insertPoint.InsertXML($newParagraphsXml)

# Move the "_GoBack" bookmark off the original closing paragraph: the
# collection resolves by name to the first (oldest) match, which is still
# the stale one on "...training data.", so deleting it leaves only the
# copy we just inserted on the "Cross v without new features" paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
